$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 3: "aantal uur" 2 -> 3, and taller row height
$ws.Range("B3").Value2 = 3
$ws.Rows.Item(3).RowHeight = 63.75

# Add a new row 4 with the same formatting as row 3
$ws.Range("A3:E3").Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A4").Value2 = 45185
$ws.Range("B4").Value2 = 2
$ws.Range("C4").Value2 = "meedere bestanden, twee extra locaties, nog meer aanpassingen aan de winkel"
$ws.Range("D4").Value2 = "de verhalen dictionary werd te groot, de speler moet nog steeds sommige items kunnen krijgen ookal hoeft hij niet meer naar de winkel"
$ws.Range("E4").Value2 = "verhalen naar een ander python bestand verplaatst, de speler kan kiezen om nog steeds naar de winkel te gaan"

$ws.Rows.Item(4).RowHeight = 60

$ws.Range("E4").Select()
